$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(6441.65, 2117.94, 297.79, 566.43, 497.33, 206.26, 264.15, 25.04)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$ws.Range("B9").Select()
